$d = $word.ActiveDocument

# Find the 1-based index of the paragraph "Cam ID -> random effect" (the
# last sub-bullet under "Covariate update", right before the blank
# paragraph that precedes "Raw data -> ...").
$i = 0
$targetIndex = -1
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Cam ID -> random effect*") {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    $target = $d.Paragraphs.Item($targetIndex)

    # Insert a new paragraph right after it; the new paragraph inherits
    # $target's pPr (ListParagraph style, ilvl 1 / numId 2).
    $target.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.Text = "Camera type -> detection"
}
